# Updated cryptos list on Tue Sep  5 17:55:30 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay plain text so numeric-looking values
# (e.g. "215.90", "0.06436") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.875.64"
$ws.Range("E2").Value = "  -0.19%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.643.65"

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.19%  "

# Row 5 - BNB
$ws.Range("D5").Value = "215.90"
$ws.Range("E5").Value = "  -0.09%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.5061"
$ws.Range("E6").Value = "  -0.69%  "

# Row 7 - USDC
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.22%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2586"
$ws.Range("E8").Value = "  +0.30%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06436"
$ws.Range("E9").Value = "  +1.29%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  +5.28%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.07801"
$ws.Range("E11").Value = "  +0.14%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "4.273"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13 - now WrappedliquidstakedEther2.0 (swapped with row 14)
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.870.57"
$ws.Range("E13").Value = "  +0.61%  "

# Row 14 - now WrappedEther (swapped with row 13)
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.638.85"
$ws.Range("E14").Value = "  +0.27%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.5639"
$ws.Range("E15").Value = "  +2.31%  "

# Row 16 - ShibaInu
$sub5 = [char]0x2085
$ws.Range("D16").Value = "0.0${sub5}7689"
$ws.Range("E16").Value = "  +0.29%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "63.38"
$ws.Range("E17").Value = "  -0.92%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.895.25"
$ws.Range("E18").Value = "  -0.20%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.29%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "193.41"
$ws.Range("E20").Value = "  -1.34%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "9.937"
$ws.Range("E22").Value = "  +0.30%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "6.136"
$ws.Range("E23").Value = "  +1.18%  "

# Row 24 - BinanceUSD
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.16%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "1.808"
$ws.Range("E25").Value = "  -5.18%  "

# Row 26 - Monero
$ws.Range("D26").Value = "141.55"
$ws.Range("E26").Value = "  -0.52%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.93%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "6.807"
$ws.Range("E28").Value = "  +0.61%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "15.53"
$ws.Range("E29").Value = "  -0.74%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.34%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "0.04954"
$ws.Range("E31").Value = "  +0.86%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "3.301"
$ws.Range("E32").Value = "  +1.63%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "3.238"
$ws.Range("E33").Value = "  +1.10%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "1.573"
$ws.Range("E34").Value = "  +1.95%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "2.387"
$ws.Range("E35").Value = "  +0.72%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Value = "0.9056"
$ws.Range("E36").Value = "  +0.78%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +0.29%  "

# Row 38 - Maker
$ws.Range("D38").Value = "1.133.86"
$ws.Range("E38").Value = "  +1.88%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "2.556"
$ws.Range("E39").Value = "  +0.74%  "

# Row 40 - VeChain
$ws.Range("D40").Value = "0.01568"
$ws.Range("E40").Value = "  +0.56%  "

# Row 41 - PaxDollar
$ws.Range("D41").Value = "0.9981"
$ws.Range("E41").Value = "  -0.29%  "

# Row 42 - FraxShare
$ws.Range("E42").Value = "  -1.78%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "0.8083"
$ws.Range("E43").Value = "  +1.63%  "

# Row 44 - Quant
$ws.Range("D44").Value = "98.83"
$ws.Range("E44").Value = "  +1.36%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.780.03"
$ws.Range("E45").Value = "  +0.55%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  -7.55%  "

# Row 47 - Aave
$ws.Range("D47").Value = "55.78"
$ws.Range("E47").Value = "  +1.77%  "

# Row 48 - Mantle
$ws.Range("D48").Value = "0.4288"
$ws.Range("E48").Value = "  -3.64%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "7.758"
$ws.Range("E49").Value = "  +2.33%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.05045"
$ws.Range("E50").Value = "  -1.77%  "

# Row 51 - Frax
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  -0.08%  "
